# tabelasProjeto.xlsx - add a "userlevel" column to the USR table, between
# "password" and "address".
#
# The USR table header lives in row 6 (A6:H6 before the edit):
#   ID | firstname | surname | e-mail | password | address | phone | cellphone
# A banner row (row 5, merged A5:H5) sits above it.
#
# We insert a brand-new column at F (pushing address/phone/cellphone one
# column to the right, extending the row-5 merge and the used range), then
# fill in the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new column before column F. This shifts the existing
# F:H columns (address, phone, cellphone) to G:I, widens the A5:H5 merged
# banner to A5:I5, and extends the sheet's used range to A3:I6 - all
# exactly like Excel's native "Insert Sheet Columns" command.
$ws.Columns.Item(6).Insert()

# New header cell for the inserted column, styled the same as its
# neighbours since Insert() carries the surrounding column formatting.
$ws.Range("F6").Value = "userlevel "

# Leave the selection where the author apparently left off.
[void]$ws.Range("F18").Select()
